# ---------------------------------------------------------------------------
# Applies the "[ADDITIONAL SCRAPING]" commit:
#   1. Inserts a new first sheet "Player Info" with ID/NAME/BATTING_HAND/
#      BOWL_STYLE columns for player 6058 (Mohammad Wasim).
#   2. Renames the MATCH_CARD_LINK column to MATCH_CODE on both the
#      "ODI Batting" and "ODI Bowling" sheets, replacing the full
#      howstat.com scorecard URL with just the numeric match code that was
#      embedded in it (?MatchCode=NNNN -> NNNN).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- helper: write a value into a cell while forcing it to stay a text
# cell (the workbook otherwise auto-coerces plain digit strings such as
# "6058" or "4564" into numbers). We flip the cell to text format, assign
# the string, then reset the style back to Normal/General so the cell is
# left with no special numeric formatting applied (matching the rest of
# the sheet's untouched data cells). ---------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before "ODI Batting" (becomes the
#    first sheet in the workbook).
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# NOTE: inserting a sheet "before" $battingSheet repoints that handle at the
# freshly inserted sheet instead of following "ODI Batting" to its new
# position, so both references must be re-fetched by name afterwards.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold/bordered/centered header styling already used by the
# other sheets' header rows by copying it over from one of them.
$battingSheet.Range("A1").Copy() | Out-Null
$playerInfo.Range("A1:D1").PasteSpecial(-4122) | Out-Null ; # xlPasteFormats

Set-TextValue $playerInfo.Range("A2") "6058"
$playerInfo.Range("B2").Value = "Mohammad Wasim"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# Pulls the numeric "MatchCode=NNNN" query parameter out of a howstat.com
# scorecard URL. ($range.Value comes back as a bare COM variant that
# doesn't behave like a normal string, so read .Value2 instead.)
function Get-MatchCode($range) {
    $url = $range.Value2
    if ($url -match 'MatchCode=(\d+)') {
        return $matches[1]
    }
    return $url
}

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> code.
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$lastRow = $battingSheet.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $battingSheet.Range("D$row")
    $code = Get-MatchCode $cell
    Set-TextValue $cell $code
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> code.
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$lastRow = $bowlingSheet.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $bowlingSheet.Range("B$row")
    $code = Get-MatchCode $cell
    Set-TextValue $cell $code
}

Write-Host "Done. Sheets now:"
foreach ($ws in $wb.Worksheets) {
    Write-Host (" - " + $ws.Name)
}
